# Updated cryptos list (prices / 1h volume change) per upstream diff.
# Values that look purely numeric (e.g. "1.00", "586.67") are entered with a
# leading apostrophe so Excel stores them as text (matching the original
# inline-string cell type) instead of silently converting them to numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.299.18"
$ws.Range("E2").Value = "  +0.58%  "
$ws.Range("D3").Value = "3.493.63"
$ws.Range("E3").Value = "  -0.15%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "'586.67"
$ws.Range("E5").Value = "  +0.31%  "
$ws.Range("D6").Value = "'134.12"
$ws.Range("E6").Value = "  +1.39%  "
$ws.Range("D7").Value = "3.494.25"
$ws.Range("E7").Value = "  -0.12%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("E9").Value = "  -1.12%  "
$ws.Range("E10").Value = "  -0.06%  "
$ws.Range("D11").Value = "'7.18"
$ws.Range("E11").Value = "  +1.02%  "
$ws.Range("E12").Value = "  -2.20%  "
$ws.Range("D13").Value = "4.088.78"
$ws.Range("E13").Value = "  -0.21%  "
$ws.Range("E14").Value = "  +1.67%  "
$ws.Range("E15").Value = "  +0.62%  "
$ws.Range("D16").Value = "3.491.95"
$ws.Range("E16").Value = "  -0.32%  "
$ws.Range("D17").Value = "64.338.10"
$ws.Range("E17").Value = "  +0.48%  "
$ws.Range("E18").Value = "  -9.24%  "
$ws.Range("D19").Value = "'9.95"
$ws.Range("E19").Value = "  -0.43%  "
$ws.Range("E20").Value = "  +1.28%  "
$ws.Range("D21").Value = "'13.60"
$ws.Range("E21").Value = "  -6.38%  "
$ws.Range("D22").Value = "'387.54"
$ws.Range("E22").Value = "  -0.95%  "
$ws.Range("D23").Value = "'0.566"
$ws.Range("E23").Value = "  -2.17%  "
$ws.Range("D24").Value = "3.633.68"
$ws.Range("E24").Value = "  -0.15%  "
$ws.Range("D25").Value = "'74.45"
$ws.Range("E25").Value = "  +1.89%  "
$ws.Range("E26").Value = "  -0.09%  "
$ws.Range("E27").Value = "  -0.82%  "
$ws.Range("E28").Value = "  +0.03%  "
$ws.Range("D29").Value = "'1.55"
$ws.Range("E29").Value = "  -2.27%  "
$ws.Range("E30").Value = "  +0.04%  "
$ws.Range("E31").Value = "  -1.42%  "
$ws.Range("E32").Value = "  +0.47%  "
$ws.Range("E33").Value = "  -0.94%  "
$ws.Range("D34").Value = "3.515.90"
$ws.Range("E34").Value = "  +0.33%  "
$ws.Range("D36").Value = "'0.149"
$ws.Range("E36").Value = "  +2.80%  "
$ws.Range("D37").Value = "'23.47"
$ws.Range("E37").Value = "  -1.66%  "
$ws.Range("E38").Value = "  -0.83%  "
$ws.Range("D39").Value = "'6.86"
$ws.Range("E39").Value = "  -1.85%  "
$ws.Range("E40").Value = "  -1.65%  "
$ws.Range("D41").Value = "'161.83"
$ws.Range("E41").Value = "  -3.50%  "
$ws.Range("E42").Value = "  -3.47%  "
$ws.Range("E43").Value = "  -1.04%  "
$ws.Range("E44").Value = "  -0.02%  "
$ws.Range("D45").Value = "'25.32"
$ws.Range("E45").Value = "  -6.44%  "
$ws.Range("D46").Value = "'41.90"
$ws.Range("E46").Value = "  -0.15%  "
$ws.Range("E47").Value = "  +0.28%  "
$ws.Range("E48").Value = "  -2.13%  "
$ws.Range("D49").Value = "'1.65"
$ws.Range("E49").Value = "  +0.66%  "
$ws.Range("D50").Value = "2.472.86"
$ws.Range("E50").Value = "  +1.42%  "
$ws.Range("D51").Value = "'6.75"
$ws.Range("E51").Value = "  -2.12%  "
